$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Insert a brand-new "2022-Q1" sheet right after "2021-Q4" (and
#    before "总计"), laid out the same way as the other per-quarter
#    fund-holding sheets (e.g. "2021-Q4").
# ------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $refSheet)
$newSheet.Name = "2022-Q1"

# Copy header-row formatting (bold/border/center style) and the
# column-A row-index style from the reference sheet so the new sheet
# matches the existing visual style exactly.
$refSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$refSheet.Range("A2:A3").Copy()
$newSheet.Range("A2:A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row
$newSheet.Cells.Item(1, 2).Value2 = "基金代码"
$newSheet.Cells.Item(1, 3).Value2 = "基金名称"
$newSheet.Cells.Item(1, 4).Value2 = "基金规模"
$newSheet.Cells.Item(1, 5).Value2 = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value2 = "仓位占比"
$newSheet.Cells.Item(1, 7).Value2 = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value2 = "仓位排名"

# Row 2 - 009869 嘉实产业先锋混合A
$newSheet.Cells.Item(2, 1).Value2 = 0
$newSheet.Cells.Item(2, 2).Value2 = "'009869"
$newSheet.Cells.Item(2, 2).Style = "Normal"
$newSheet.Cells.Item(2, 3).Value2 = "'嘉实产业先锋混合A"
$newSheet.Cells.Item(2, 3).Style = "Normal"
$newSheet.Cells.Item(2, 4).Value2 = "'17.16"
$newSheet.Cells.Item(2, 4).Style = "Normal"
$newSheet.Cells.Item(2, 5).Value2 = "'89.91"
$newSheet.Cells.Item(2, 5).Style = "Normal"
$newSheet.Cells.Item(2, 6).Value2 = "'4.10"
$newSheet.Cells.Item(2, 6).Style = "Normal"
$newSheet.Cells.Item(2, 7).Value2 = "'0.7036"
$newSheet.Cells.Item(2, 7).Style = "Normal"
$newSheet.Cells.Item(2, 8).Value2 = 7

# Row 3 - 009870 嘉实产业先锋混合C
$newSheet.Cells.Item(3, 1).Value2 = 1
$newSheet.Cells.Item(3, 2).Value2 = "'009870"
$newSheet.Cells.Item(3, 2).Style = "Normal"
$newSheet.Cells.Item(3, 3).Value2 = "'嘉实产业先锋混合C"
$newSheet.Cells.Item(3, 3).Style = "Normal"
$newSheet.Cells.Item(3, 4).Value2 = "'3.71"
$newSheet.Cells.Item(3, 4).Style = "Normal"
$newSheet.Cells.Item(3, 5).Value2 = "'89.91"
$newSheet.Cells.Item(3, 5).Style = "Normal"
$newSheet.Cells.Item(3, 6).Value2 = "'4.10"
$newSheet.Cells.Item(3, 6).Style = "Normal"
$newSheet.Cells.Item(3, 7).Value2 = "'0.1521"
$newSheet.Cells.Item(3, 7).Style = "Normal"
$newSheet.Cells.Item(3, 8).Value2 = 7

# ------------------------------------------------------------------
# 2) Update the "总计" (totals) sheet: push its existing rows down by
#    one and insert a new "2022-Q1" summary row at the top of the data
#    (row 2).
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

# Extend the row-index column style down into the new last row (7)
# before shifting values into it.
$totalSheet.Range("A6").Copy()
$totalSheet.Range("A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Shift rows 2..6 down to 3..7, bottom-up so sources aren't clobbered.
for ($r = 6; $r -ge 2; $r--) {
    $b = $totalSheet.Cells.Item($r, 2).Value2
    $c = $totalSheet.Cells.Item($r, 3).Value2
    $d = $totalSheet.Cells.Item($r, 4).Value2
    $totalSheet.Cells.Item($r + 1, 1).Value2 = $r - 1
    $totalSheet.Cells.Item($r + 1, 2).Value2 = $b
    $totalSheet.Cells.Item($r + 1, 3).Value2 = $c
    $totalSheet.Cells.Item($r + 1, 4).Value2 = $d
}

# New top data row: 2022-Q1 summary.
$totalSheet.Cells.Item(2, 1).Value2 = 0
$totalSheet.Cells.Item(2, 2).Value2 = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value2 = 2
$totalSheet.Cells.Item(2, 4).Value2 = 0.86
